# Update the CPD receiver-operating-characteristic summary sheet:
# - split each State/Non-State/One-Sided metric column into separate mean/std columns
# - rename the CART algorithm row to DTREE
# - drop the NB (Naive Bayes) row, shrinking the data from 8 rows to 7
# - refresh all numeric results to the new values

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the last data row (old NB row, row 9) since the dataset now has 7 rows instead of 8
$ws.Rows.Item(9).Delete() | Out-Null

# Update header row (B1:E1 text changes, add F1:H1)
$ws.Range("B1").Value = "Algorithm"
$ws.Range("C1").Value = "State Based mean"
$ws.Range("D1").Value = "State Based std"
$ws.Range("E1").Value = "Non State mean"
$ws.Range("F1").Value = "Non State std"
$ws.Range("G1").Value = "One Sided mean"
$ws.Range("H1").Value = "One Sided std"

# Apply the existing header style (bold, centered, bordered) to the new F1:H1 header cells
# by copying the format from B1 (re-uses the existing style instead of creating a new one)
$ws.Range("B1").Copy() | Out-Null
$ws.Range("F1:H1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Update data rows: B (algorithm label) stays/renamed, C/D become State Based mean/std,
# E/F become Non State mean/std, G/H become One Sided mean/std
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "LR"
$ws.Range("C2").Value = 0.9006795195954487
$ws.Range("D2").Value = 0.03151325806880326
$ws.Range("E2").Value = 0.8759856630824372
$ws.Range("F2").Value = 0.04952076963298688
$ws.Range("G2").Value = 0.8870141784820683
$ws.Range("H2").Value = 0.03257053301693878

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "LDA"
$ws.Range("C3").Value = 0.9059813527180784
$ws.Range("D3").Value = 0.03189442994976523
$ws.Range("E3").Value = 0.9114183307731695
$ws.Range("F3").Value = 0.0231411056431009
$ws.Range("G3").Value = 0.9060884070058381
$ws.Range("H3").Value = 0.02166676614662823

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "KNN"
$ws.Range("C4").Value = 0.9396886852085966
$ws.Range("D4").Value = 0.0219441420580556
$ws.Range("E4").Value = 0.9436251920122889
$ws.Range("F4").Value = 0.02196728000553226
$ws.Range("G4").Value = 0.950767306088407
$ws.Range("H4").Value = 0.01919342084510206

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "DTREE"
$ws.Range("C5").Value = 0.747115992414665
$ws.Range("D5").Value = 0.03050516544368534
$ws.Range("E5").Value = 0.6458269329237072
$ws.Range("F5").Value = 0.05878843102410158
$ws.Range("G5").Value = 0.6973227689741451
$ws.Range("H5").Value = 0.04198481601140375

$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "RTREE"
$ws.Range("C6").Value = 0.8811551833122628
$ws.Range("D6").Value = 0.03080224056728806
$ws.Range("E6").Value = 0.8598054275473629
$ws.Range("F6").Value = 0.05067006027591416
$ws.Range("G6").Value = 0.8997497914929108
$ws.Range("H6").Value = 0.02294896883659027

$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "XTREE"
$ws.Range("C7").Value = 0.8961994310998735
$ws.Range("D7").Value = 0.02565312876135361
$ws.Range("E7").Value = 0.8791858678955453
$ws.Range("F7").Value = 0.0410558734383163
$ws.Range("G7").Value = 0.879674728940784
$ws.Range("H7").Value = 0.0176100125399712

$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "SVM"
$ws.Range("C8").Value = 0.8988701011378002
$ws.Range("D8").Value = 0.03434592732462118
$ws.Range("E8").Value = 0.8630568356374807
$ws.Range("F8").Value = 0.04357357507581917
$ws.Range("G8").Value = 0.89883236030025
$ws.Range("H8").Value = 0.01964904255532755

